$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 4.3
$ws.Range("M3").Value = 3.8
$ws.Range("M4").Value = 3.6
$ws.Range("M5").Value = 3.4
$ws.Range("M6").Value = 2.8
$ws.Range("I7").Value = 0.112906918401084
$ws.Range("M7").Value = 2.5
$ws.Range("M8").Value = 2
